$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3907
$ws.Range("I3").Value = 7489
$ws.Range("K3").Value = 3963
$ws.Range("I4").Value = 1795
$ws.Range("J4").Value = 1820
$ws.Range("K4").Value = 806
$ws.Range("K5").Value = 276
$ws.Range("K6").Value = 4504
$ws.Range("J7").Value = 29289
$ws.Range("K7").Value = 13456

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 262
$ws.Range("K3").Value = 272
$ws.Range("K4").Value = 51
$ws.Range("K6").Value = 308
$ws.Range("K7").Value = 915

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 95
$ws.Range("K7").Value = 283

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K6").Value = 161
$ws.Range("K7").Value = 555

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 226

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 93
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 115
$ws.Range("K4").Value = 46
$ws.Range("K7").Value = 390
$ws.Range("K8").Value = 915
$ws.Range("K9").Value = 55
$ws.Range("K10").Value = 77
$ws.Range("K19").Value = 413
$ws.Range("K20").Value = 299
$ws.Range("K23").Value = 135
$ws.Range("K27").Value = 132
$ws.Range("K29").Value = 710
$ws.Range("K31").Value = 146
$ws.Range("K33").Value = 555
$ws.Range("K34").Value = 66
$ws.Range("K35").Value = 20
$ws.Range("K36").Value = 169
$ws.Range("K40").Value = 34
$ws.Range("K42").Value = 475
$ws.Range("K44").Value = 125
$ws.Range("K47").Value = 76
$ws.Range("K48").Value = 177
$ws.Range("K50").Value = 74
$ws.Range("K51").Value = 162
$ws.Range("K52").Value = 368
$ws.Range("K53").Value = 182
$ws.Range("K55").Value = 151
$ws.Range("K60").Value = 87
$ws.Range("J63").Value = 104
$ws.Range("K63").Value = 45
$ws.Range("K65").Value = 315
$ws.Range("K67").Value = 523
$ws.Range("K72").Value = 65
$ws.Range("K76").Value = 191
$ws.Range("K78").Value = 165
$ws.Range("K79").Value = 349
$ws.Range("K83").Value = 283
$ws.Range("K85").Value = 605
$ws.Range("K86").Value = 91
$ws.Range("K87").Value = 18
$ws.Range("K88").Value = 154
$ws.Range("K89").Value = 188
$ws.Range("K91").Value = 146
$ws.Range("K94").Value = 166
$ws.Range("K95").Value = 226
$ws.Range("K97").Value = 115
$ws.Range("K99").Value = 236
$ws.Range("J101").Value = 29289
$ws.Range("K101").Value = 13456

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 49
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 155
$ws.Range("K3").Value = 180
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 523

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 204
$ws.Range("K3").Value = 248
$ws.Range("K6").Value = 201
$ws.Range("K7").Value = 710

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 41
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 124
$ws.Range("K5").Value = 15
$ws.Range("K6").Value = 126
$ws.Range("K7").Value = 413

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 129
$ws.Range("K3").Value = 155
$ws.Range("K7").Value = 475

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 48
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 41
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 36
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K5").Value = 11
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 103
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 299

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 66
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 169

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 123
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 390

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 47
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K2").Value = 1
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 16
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 35
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 48
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 44
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 217
$ws.Range("K6").Value = 139
$ws.Range("K7").Value = 605

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 8
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 95
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 368

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 18
